$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full content of row 5 and row 6 (everything except the
# handful of columns that already held identical values in both rows, such
# as D/S/T/U/V/W/Y/AA/AD/AE/AG). Only touch the cells that actually change,
# to avoid disturbing already-matching columns (e.g. the "2026-02-07" text
# in Y/AA, which must stay plain text rather than being re-interpreted as a
# date by a round-trip through .Value2).

# New values for row 5 (taken from row 6's original content).
$row5set = @{
    "A"  = 131067473
    "B"  = 57884
    "E"  = 100109
    "F"  = "Tretåig hackspett"
    "G"  = "Picoides tridactylus"
    "H"  = "(Linnaeus, 1758)"
    "M"  = "färska spår"
    "P"  = "Åbogen, Jmt"
    "Q"  = 465809
    "R"  = 7046259
    "Z"  = "15:46"
    "AB" = "15:46"
    "AC" = "Färska ringhack"
    "AW" = "Elin Albrechtsson"
    "AX" = "Elin Albrechtsson"
}
$row5clear = @("J", "L", "AF", "AH")

# New values for row 6 (taken from row 5's original content).
$row6set = @{
    "A"  = 131067826
    "B"  = 79243
    "E"  = 6425
    "F"  = "Garnlav"
    "G"  = "Alectoria sarmentosa"
    "H"  = "(Ach.) Ach."
    "P"  = "Långan Öst, Jmt"
    "Q"  = 465891
    "R"  = 7046290
    "AH" = "Granskog"
    "AW" = "Kristian Zackrisson"
    "AX" = "Kristian Zackrisson"
}
$row6clear = @("J", "L", "M", "Z", "AB", "AC")

foreach ($col in $row5set.Keys) {
    $ws.Range($col + "5").Value2 = $row5set[$col]
}
foreach ($col in $row5clear) {
    $ws.Range($col + "5").Value2 = $null
}

foreach ($col in $row6set.Keys) {
    $ws.Range($col + "6").Value2 = $row6set[$col]
}
foreach ($col in $row6clear) {
    $ws.Range($col + "6").Value2 = $null
}
